$d = $word.ActiveDocument

# --- Fix 1: paragraph with the YouTube hyperlink loses <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr> from its <w:pPr> ---
$hlPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*UzPjHQX5-9A*") {
        $hlPara = $p
        break
    }
}
if ($hlPara -eq $null) { throw "hyperlink paragraph not found" }

# Remove the paragraph's own text/hyperlink content, keep the (old) paragraph mark with its pPr for now
$hlContentRange = $d.Range($hlPara.Range.Start, $hlPara.Range.End - 1)
$hlContentRange.Delete()

# Re-fetch the (now empty) paragraph and insert the replacement paragraph (correct pPr, no rFonts hint)
# followed by a throw-away empty paragraph which will absorb the merge onto the OLD (hint-bearing) pPr.
$hlParaEmpty = $d.Paragraphs.Item($hlPara.Index)
$insAt = $d.Range($hlParaEmpty.Range.Start, $hlParaEmpty.Range.Start)
$hlFixXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:ind w:left="360" w:firstLineChars="0" w:firstLine="0"/></w:pPr><w:hyperlink r:id="rId8" w:history="1"><w:r><w:t>https://www.youtube.com/watch?v=UzPjHQX5-9A</w:t></w:r></w:hyperlink></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insAt.InsertXML($hlFixXml)

# Re-apply the Hyperlink character style to the URL run (style refs aren't resolved from the minimal XML package)
$styleRange = $d.Content.Duplicate
$styleRange.Find.Execute("https://www.youtube.com/watch?v=UzPjHQX5-9A")
$styleRange.Style = $d.Styles.Item("a7")

# Clean up the two leftover empty paragraphs produced by the merge trick (dummy + old hint-bearing one)
$newHlPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*UzPjHQX5-9A*") {
        $newHlPara = $p
        break
    }
}
$afterA = $d.Paragraphs.Item($newHlPara.Index + 1)
$afterB = $d.Paragraphs.Item($newHlPara.Index + 2)
$cleanupR = $d.Range($afterA.Range.Start, $afterB.Range.End)
$cleanupR.Delete()
$afterA2 = $d.Paragraphs.Item($newHlPara.Index + 1)
$afterB2 = $d.Paragraphs.Item($newHlPara.Index + 2)
$cleanupR2 = $d.Range($afterA2.Range.Start, $afterB2.Range.Start)
$cleanupR2.Delete()

# --- Fix 2: replace the final (empty) paragraph with the new ORB / Tracking-Mapping-Relocation content ---
$lastPara = $d.Paragraphs.Last
$insAt2 = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$contentXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:ind w:firstLine="360"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">ORB features : </w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:t>ORB: an efficient alternative to SIFT or SURF</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/></w:pPr><w:r><w:t>ORB = oFAST + rBRIEF</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/></w:pPr><w:r><w:t>oFAST: FAST Keypoint Orientation</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/></w:pPr><w:r><w:t>rBRIEF: Rotation-Aware Brief</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:t>Brief: Binary robust independent elementary features.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:ind w:firstLine="360"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Tracking, Mapping, Relocation</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Tracking : a valid ORB binary descriptor</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">Mapping: </w:t></w:r><w:r><w:t>ORB features are used both for mapping, and for the place</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>recognition. Place recognition combines a Bag of Words built from the ORB</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>binary descriptors, with the covisibility graph that determines all the keyframes</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>that are observing the same 3D scene region</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="360"/></w:pPr><w:r><w:t>Covisilibilty Graph 的顶点是相机的Pose，而边是Pose-Pose的变换关系——所以也算是Pose Graph 一种吧。当两个相机看到相似的空间点时，它们对应的Pose就会产生联系（我们就可以根据这些空间点在照片上的投影计算两个相机间的运动）。根据观测到的空间点的数量，给这个边加上一个权值，度量这个边的可信程度。</w:t></w:r><w:r><w:t>Covisibility Graph是一个无向有权图(graph),这个概念最早来自2010的文章[Closing Loops Without Places]。</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>简单来说，每个</w:t></w:r><w:r><w:t>node就是关键帧，edge的权重就是两个关键帧找到足够多的相同的 3d 点的数目。</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insAt2.InsertXML($contentXml)

# Clean up the trailing empty paragraph left over from the merge trick
$p1 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$cleanupR3 = $d.Range($p1.Range.End - 1, $p2.Range.End)
$cleanupR3.Delete()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
Write-Host "Last paragraph text:" $d.Paragraphs.Last.Range.Text
